$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.129390239715576
$ws.Range("B1").Value = 2.47679877281189
$ws.Range("C1").Value = 1.821427583694458
$ws.Range("D1").Value = 1.692352294921875
$ws.Range("E1").Value = 1.592419266700745
